$wb = $excel.ActiveWorkbook

# The previously-active sheet ("TUC - WWKYS Submenu Names") is no longer the
# selected tab; reset its lingering selection back to the top-left cell.
$wsWwkys = $wb.Worksheets.Item("TUC - WWKYS Submenu Names")
[void]$wsWwkys.Range("A1").Select()

# Add "TUC - FDTOC Submenu Names" after the current last sheet
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsFdtoc = $wb.Worksheets.Add($null, $lastSheet)
$wsFdtoc.Name = "TUC - FDTOC Submenu Names"

# Add "TUC - FWYNTK Submenu Names" after "TUC - FDTOC Submenu Names"
$wsFwyntk = $wb.Worksheets.Add($null, $wsFdtoc)
$wsFwyntk.Name = "TUC - FWYNTK Submenu Names"

# Fill content for "TUC - FDTOC Submenu Names"
$wsFdtoc.Range("A1").Value = "Weekly Update from CEO Ed Bastian"
$wsFdtoc.Range("A2").Value = "How Delta is Supporting Medical Volunteers"

# Fill content for "TUC - FWYNTK Submenu Names"
$wsFwyntk.Range("A1").Value = "Coronavirus Regional Restrictions"
$wsFwyntk.Range("A2").Value = "Delta Temporarily Closes Select Airports"
$wsFwyntk.Range("A3").Value = "Things to Know When You Travel with a Partner Airline"

# Make the new last sheet ("TUC - FWYNTK Submenu Names") the active/selected one
$wsFwyntk.Activate()
